$wb = $excel.ActiveWorkbook
$wsSources = $wb.Worksheets.Item("Sources")
$wsClaims = $wb.Worksheets.Item("Claims")
$wsFigures = $wb.Worksheets.Item("Figures")

# ---- Sources sheet ----
# Row 90
$wsSources.Range("D90").Value = "_figures/figures_data.xlsx"
$wsSources.Range("I90").Value = "[AUDIT 2026-02-08] Canonical figure-data workbook migrated to _figures/figures_data.xlsx."

# Row 105
$wsSources.Range("D105").Value = "sources/datasets/Whitepaper_Master_Data.xlsx"
$wsSources.Range("I105").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for Grand View figures used in report models."

# Row 106
$wsSources.Range("D106").Value = "sources/datasets/Pet_Ownership_Market_Data_Sourced.xlsx"
$wsSources.Range("I106").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for Euromonitor pet-care values used in report models."

# Row 107
$wsSources.Range("D107").Value = "sources/datasets/Nutraceuticals_Delivery_Formats_Data.xlsx"
$wsSources.Range("I107").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for NBJ format/category data used in report models."

# Row 108
$wsSources.Range("D108").Value = "sources/datasets/Livestock_Segment_Complete_Data_Sources.xlsx"
$wsSources.Range("I108").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for Future Market Insights feed-additive estimates used in report models."

# Row 109
$wsSources.Range("D109").Value = "sources/datasets/Figure5_Probiotics_Share.csv"
$wsSources.Range("I109").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for MarketsandMarkets probiotics split used in report models."

# Row 110
$wsSources.Range("D110").Value = "sources/academic/FEDIAF-Facts-Figures-2025.pdf"

# Row 111
$wsSources.Range("D111").Value = "sources/datasets/Figure1_Pet_Ownership.csv"
$wsSources.Range("I111").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for APPA ownership metrics used in figures."

# Row 112
$wsSources.Range("D112").Value = "sources/datasets/Figure9_Livestock_Trends.csv"
$wsSources.Range("I112").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for FAO SOFIA livestock/aquaculture trend data used in figures."

# Row 113
$wsSources.Range("D113").Value = "sources/datasets/Livestock_Breeding_EU_vs_USA_Data.xlsx"
$wsSources.Range("I113").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for Eurostat livestock/meat comparison used in report models."

# Row 114
$wsSources.Range("D114").Value = "sources/datasets/Nutraceuticals_Functional_Segments_Data.xlsx"
$wsSources.Range("I114").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Derived dataset proxy for Mordor Intelligence market segmentation used in report models."

# Row 115
$wsSources.Range("D115").Value = "sources/academic/Nutraceuticals, Social Interaction, and Psychophysiological Influence on Pet Health and Well-Being- Focus on Dogs and Cats.pdf"

# Row 116
$wsSources.Range("D116").Value = "sources/reports/zoetis_2024_annual_report.pdf"
$wsSources.Range("I116").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Primary filing located in archived inputs and copied into sources/reports."

# Row 117
$wsSources.Range("D117").Value = "sources/internal/20260115_VC_PE_Portfolio.xlsx"
$wsSources.Range("I117").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Internal investor compilation used as in-repo proxy; original Elanco deck not found by filename."

# Row 118
$wsSources.Range("D118").Value = "sources/reports/dsm_firmenich_2024_integrated_annual_report.pdf"
$wsSources.Range("I118").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Primary filing located in archived inputs and copied into sources/reports."

# Row 119
$wsSources.Range("D119").Value = "sources/reports/swedencare_annual_report_2024.pdf"
$wsSources.Range("I119").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Primary filing located in archived inputs and copied into sources/reports."

# Row 120
$wsSources.Range("D120").Value = "sources/reports/virbac_annual_report_2024.pdf"
$wsSources.Range("I120").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Primary filing located in archived inputs and copied into sources/reports."

# Row 121
$wsSources.Range("D121").Value = "sources/reports/dechra_annual_report_2023.pdf"
$wsSources.Range("I121").Value = "Y [AUDIT 2026-02-08] SOURCE FILE NOT FOUND IN REPO; author must attach or replace with an on-disk file path. [AUDIT 2026-02-08 SOURCE-RECOVERY] Primary filing located in archived inputs and copied into sources/reports."

# ---- Claims sheet ----
# Row 4
$wsClaims.Range("D4").Value = "AUTHOR-CHECK"
$wsClaims.Range("F4").Value = "N"
$wsClaims.Range("J4").Value = "Figure 4 Regional Market sums to 6.0B. Valuation range based on M&A transaction analysis [AUDIT 2026-02-08] Needs external support for valuation multiple (15-20x EBITDA). Candidate filings: S118/S119/S120."

# Row 6
$wsClaims.Range("D6").Value = "AUTHOR-CHECK"
$wsClaims.Range("F6").Value = "N"
$wsClaims.Range("J6").Value = "Derived from 13-6. S039 has 8B for Livestock Premix. Valuation range based on livestock sector M&A comps [AUDIT 2026-02-08] Market size is calculation-backed; valuation multiple requires external support (S118/S119/S120)."

# Row 8
$wsClaims.Range("D8").Value = "AUTHOR-CHECK"
$wsClaims.Range("F8").Value = "N"
$wsClaims.Range("J8").Value = "Source pending. DTC vs retail margin analysis [AUDIT 2026-02-08] Retail vs DTC margin compression requires external channel data source."

# Row 10
$wsClaims.Range("D10").Value = "S085"
$wsClaims.Range("F10").Value = "N"
$wsClaims.Range("J10").Value = "Matches Table_US_vs_EU.csv. Pharma-grade vs commodity pricing analysis [AUDIT 2026-02-08] Mapped to Table_US_vs_EU dataset."

# Row 11
$wsClaims.Range("D11").Value = "S086"
$wsClaims.Range("F11").Value = "N"
$wsClaims.Range("J11").Value = "Matches Timeline_Regulations.csv. Veterinary channel CLV impact modeling [AUDIT 2026-02-08] Mapped to Timeline_Regulations dataset."

# Row 12
$wsClaims.Range("D12").Value = "S086"
$wsClaims.Range("F12").Value = "N"
$wsClaims.Range("J12").Value = "Matches Timeline_Regulations.csv. Clinical trial cost estimates from industry practice [AUDIT 2026-02-08] Mapped to Timeline_Regulations dataset."

# Row 40
$wsClaims.Range("D40").Value = "AUTHOR-CHECK"
$wsClaims.Range("F40").Value = "N"
$wsClaims.Range("J40").Value = "IP holders margin claim. IP holder margin profile from value chain modeling [AUDIT 2026-02-08] Margin band needs externally citable support."

# Row 41
$wsClaims.Range("D41").Value = "AUTHOR-CHECK"
$wsClaims.Range("F41").Value = "N"
$wsClaims.Range("J41").Value = "Commodity margin claim. Commodity supplier margins from industry benchmarks [AUDIT 2026-02-08] Margin band needs externally citable support."

# Row 42
$wsClaims.Range("D42").Value = "AUTHOR-CHECK"
$wsClaims.Range("F42").Value = "N"
$wsClaims.Range("J42").Value = "CDMO market penetration. CDMO penetration estimate [AUDIT 2026-02-08] CDMO penetration claim requires source."

# Row 43
$wsClaims.Range("D43").Value = "AUTHOR-CHECK"
$wsClaims.Range("F43").Value = "N"
$wsClaims.Range("J43").Value = "CDMO margin claim. CDMO margin analysis [AUDIT 2026-02-08] CDMO margin claim requires source."

# Row 44
$wsClaims.Range("D44").Value = "AUTHOR-CHECK"
$wsClaims.Range("F44").Value = "N"
$wsClaims.Range("J44").Value = "DTC margin after CAC. DTC margin calc after CAC [AUDIT 2026-02-08] DTC margin claim requires source."

# Row 45
$wsClaims.Range("D45").Value = "AUTHOR-CHECK"
$wsClaims.Range("F45").Value = "N"
$wsClaims.Range("J45").Value = "R&D cost claim. Pharma R&D cost benchmarks [AUDIT 2026-02-08] Molecule development cost claim requires source."

# Row 47
$wsClaims.Range("D47").Value = "AUTHOR-CHECK"
$wsClaims.Range("F47").Value = "N"
$wsClaims.Range("J47").Value = "Valuation multiple for commodity players. Commodity tier multiples from comp analysis [AUDIT 2026-02-08] Tier-2 multiple requires source support."

# Row 22
$wsClaims.Range("D22").Value = "AUTHOR-CHECK"
$wsClaims.Range("J22").Value = "Segment size mismatch with Master Excel [AUDIT 2026-02-08] Text value conflicts with S089 Figure 18; reconcile before verification."

# Row 23
$wsClaims.Range("D23").Value = "AUTHOR-CHECK"
$wsClaims.Range("J23").Value = "Segment size mismatch with Master Excel [AUDIT 2026-02-08] Text value conflicts with S089 Figure 18; reconcile before verification."

# Row 24
$wsClaims.Range("D24").Value = "AUTHOR-CHECK"
$wsClaims.Range("J24").Value = "Segment size mismatch with Master Excel [AUDIT 2026-02-08] Text value conflicts with S089 Figure 18; reconcile before verification."

# Row 25
$wsClaims.Range("D25").Value = "AUTHOR-CHECK"
$wsClaims.Range("J25").Value = "Segment size mismatch with Master Excel [AUDIT 2026-02-08] Text value conflicts with S089 Figure 18; reconcile before verification."

# ---- Figures sheet ----
# Clear empty notes cells in rows 2-21
$wsFigures.Range("I2").ClearContents()
$wsFigures.Range("I3").ClearContents()
$wsFigures.Range("I4").ClearContents()
$wsFigures.Range("I5").ClearContents()
$wsFigures.Range("I6").ClearContents()
$wsFigures.Range("I7").ClearContents()
$wsFigures.Range("I8").ClearContents()
$wsFigures.Range("I9").ClearContents()
$wsFigures.Range("I10").ClearContents()
$wsFigures.Range("I11").ClearContents()
$wsFigures.Range("I12").ClearContents()
$wsFigures.Range("I13").ClearContents()
$wsFigures.Range("I14").ClearContents()
$wsFigures.Range("I15").ClearContents()
$wsFigures.Range("I16").ClearContents()
$wsFigures.Range("I17").ClearContents()
$wsFigures.Range("I18").ClearContents()
$wsFigures.Range("I19").ClearContents()
$wsFigures.Range("I20").ClearContents()
$wsFigures.Range("I21").ClearContents()

# Row 22
$wsFigures.Range("E22").Value = "Figure 36"
$wsFigures.Range("I22").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 23
$wsFigures.Range("E23").Value = "Figure 5"
$wsFigures.Range("I23").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 24
$wsFigures.Range("E24").Value = "Figure 6"
$wsFigures.Range("I24").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 25
$wsFigures.Range("E25").Value = "Figure 16"
$wsFigures.Range("I25").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 26
$wsFigures.Range("E26").Value = "Figure 17"
$wsFigures.Range("I26").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 27
$wsFigures.Range("E27").Value = "Figure 7"
$wsFigures.Range("I27").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 28
$wsFigures.Range("E28").Value = "Figure 8"
$wsFigures.Range("I28").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 29
$wsFigures.Range("E29").Value = "Figure 9"
$wsFigures.Range("I29").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 30
$wsFigures.Range("E30").Value = "Figure 10"
$wsFigures.Range("I30").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 31
$wsFigures.Range("E31").Value = "Figure 11"
$wsFigures.Range("I31").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 32
$wsFigures.Range("E32").Value = "Figure 12"
$wsFigures.Range("I32").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 33
$wsFigures.Range("E33").Value = "Figure 15"
$wsFigures.Range("I33").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 34
$wsFigures.Range("E34").Value = "Figure 3"
$wsFigures.Range("I34").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 35
$wsFigures.Range("E35").Value = "Figure 45"
$wsFigures.Range("I35").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 36
$wsFigures.Range("E36").Value = "Figure 4"
$wsFigures.Range("I36").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 37
$wsFigures.Range("E37").Value = "Figure 2"
$wsFigures.Range("I37").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

# Row 38
$wsFigures.Range("E38").Value = "Figure 1"
$wsFigures.Range("I38").Value = "[AUDIT 2026-02-08] excel_tab remapped to canonical figures_data workbook."

Write-Output "done"